$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 98
$ws.Range("I38").Value = 98
$ws.Range("K38").Value = 294
$ws.Range("M38").Value = 78

$ws.Range("H39").Value = 1125.375
$ws.Range("I39").Value = 83.833336
$ws.Range("J39").Value = 4250
$ws.Range("K39").Value = 251.500008
$ws.Range("L39").Value = 12750
$ws.Range("M39").Value = 44.49999199999999
$ws.Range("N39").Value = -13342

$ws.Range("H40").Value = 1976.75
$ws.Range("I40").Value = 1833.3334
$ws.Range("J40").Value = 2062.8
$ws.Range("K40").Value = 1833.3334
$ws.Range("L40").Value = 2062.8
$ws.Range("M40").Value = -1658.3334
$ws.Range("N40").Value = -2412.8

$ws.Range("H51").Value = 96789.45
$ws.Range("I51").Value = 204999.4
$ws.Range("J51").Value = 6614.5
$ws.Range("K51").Value = 204999.4
$ws.Range("L51").Value = 6614.5
$ws.Range("M51").Value = -204515.4
$ws.Range("N51").Value = -7582.5

$ws.Range("H52").Value = 3057.625
$ws.Range("I52").Value = 3692.2
$ws.Range("J52").Value = 2000
$ws.Range("K52").Value = 11076.6
$ws.Range("L52").Value = 6000
$ws.Range("M52").Value = -10916.6
$ws.Range("N52").Value = -6320

$ws.Range("H57").Value = 97500
$ws.Range("J57").Value = 97500
$ws.Range("L57").Value = 292500
$ws.Range("N57").Value = -293498

$ws.Range("H61").Value = 1357.5555
$ws.Range("I61").Value = 286.33334
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 859.0000200000001
$ws.Range("L61").Value = 10500
$ws.Range("M61").Value = -687.0000200000001
$ws.Range("N61").Value = -10844

$ws.Range("H86").Value = 83340890
$ws.Range("I86").Value = 76930820
$ws.Range("K86").Value = 76930820
$ws.Range("M86").Value = -76929697

$ws.Range("H89").Value = 83340890
$ws.Range("I89").Value = 76930820
$ws.Range("K89").Value = 384654100
$ws.Range("M89").Value = -384648484

$ws.Range("H138").Value = 10370.556
$ws.Range("I138").Value = 30862.25
$ws.Range("J138").Value = 4515.7856
$ws.Range("K138").Value = 92586.75
$ws.Range("L138").Value = 13547.3568
$ws.Range("M138").Value = -87446.75
$ws.Range("N138").Value = -23827.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2277.4338
$ws.Range("I32").Value = 1251.7042
$ws.Range("K32").Value = 1251.7042
$ws.Range("M32").Value = -964.7041999999999

$ws.Range("H61").Value = 1092631
$ws.Range("I61").Value = 32026.223
$ws.Range("J61").Value = 3102197.8
$ws.Range("K61").Value = 32026.223
$ws.Range("L61").Value = 3102197.8
$ws.Range("M61").Value = -31814.223
$ws.Range("N61").Value = -3102621.8

$ws.Range("H122").Value = 1785.6818
$ws.Range("I122").Value = 1364.3
$ws.Range("K122").Value = 4092.9
$ws.Range("M122").Value = -1642.9

$ws.Range("H136").Value = 1092631
$ws.Range("I136").Value = 32026.223
$ws.Range("J136").Value = 3102197.8
$ws.Range("K136").Value = 96078.66900000001
$ws.Range("L136").Value = 9306593.399999999
$ws.Range("M136").Value = -93528.66900000001
$ws.Range("N136").Value = -9311693.399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 25589.5
$ws.Range("I107").Value = 27596.4
$ws.Range("K107").Value = 27596.4
$ws.Range("M107").Value = -25676.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5164.735
$ws.Range("I31").Value = 3907.3845
$ws.Range("J31").Value = 5398.2427
$ws.Range("K31").Value = 3907.3845
$ws.Range("L31").Value = 5398.2427
$ws.Range("M31").Value = -3612.3845
$ws.Range("N31").Value = -5988.2427

$ws.Range("H34").Value = 5164.735
$ws.Range("I34").Value = 3907.3845
$ws.Range("J34").Value = 5398.2427
$ws.Range("K34").Value = 3907.3845
$ws.Range("L34").Value = 5398.2427
$ws.Range("M34").Value = -3705.3845
$ws.Range("N34").Value = -5802.2427

$ws.Range("H58").Value = 2546.8125
$ws.Range("I58").Value = 2231.818
$ws.Range("K58").Value = 2231.818
$ws.Range("M58").Value = -2028.818

$ws.Range("H62").Value = 4432.273
$ws.Range("I62").Value = 3721
$ws.Range("J62").Value = 5025
$ws.Range("K62").Value = 3721
$ws.Range("L62").Value = 5025
$ws.Range("M62").Value = -3097
$ws.Range("N62").Value = -6273

$ws.Range("H65").Value = 4432.273
$ws.Range("I65").Value = 3721
$ws.Range("J65").Value = 5025
$ws.Range("K65").Value = 18605
$ws.Range("L65").Value = 25125
$ws.Range("M65").Value = -15485
$ws.Range("N65").Value = -31365

$ws.Range("H107").Value = 1713.7142
$ws.Range("I107").Value = 1900
$ws.Range("J107").Value = 1465.3334
$ws.Range("K107").Value = 1900
$ws.Range("L107").Value = 1465.3334
$ws.Range("M107").Value = 20
$ws.Range("N107").Value = -5305.3334

$ws.Range("H136").Value = 2546.8125
$ws.Range("I136").Value = 2231.818
$ws.Range("K136").Value = 6695.454000000001
$ws.Range("M136").Value = -4145.454000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 14974.667
$ws.Range("J105").Value = 16305.111
$ws.Range("L105").Value = 48915.333
$ws.Range("N105").Value = -54157.333

$ws.Range("H107").Value = 66667468
$ws.Range("I107").Value = 503.66666
$ws.Range("K107").Value = 1510.99998
$ws.Range("M107").Value = 409.0000199999999

$ws.Range("H115").Value = 713
$ws.Range("I115").Value = 784.3333
$ws.Range("J115").Value = 499
$ws.Range("K115").Value = 2352.9999
$ws.Range("L115").Value = 1497
$ws.Range("M115").Value = -1177.9999
$ws.Range("N115").Value = -3847

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""

$ws.Range("H131").Value = 4133879.8
$ws.Range("J131").Value = 1937.6364
$ws.Range("L131").Value = 5812.9092
$ws.Range("N131").Value = -15892.9092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = ""

$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = ""

$ws.Range("H62").Value = 45999
$ws.Range("I62").Value = 45999
$ws.Range("K62").Value = 45999
$ws.Range("M62").Value = -45313

$ws.Range("H65").Value = 45999
$ws.Range("I65").Value = 45999
$ws.Range("K65").Value = 137997
$ws.Range("M65").Value = -134565

$ws.Range("H74").Value = 40000
$ws.Range("J74").Value = 40000
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41872

$ws.Range("H77").Value = 40000
$ws.Range("J77").Value = 40000
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -129360

$ws.Range("H80").Value = 20065184
$ws.Range("J80").Value = 45573530
$ws.Range("L80").Value = 45573530
$ws.Range("N80").Value = -45575526

$ws.Range("H83").Value = 20065184
$ws.Range("J83").Value = 45573530
$ws.Range("L83").Value = 227867650
$ws.Range("N83").Value = -227877634

$ws.Range("H92").Value = 49999
$ws.Range("J92").Value = 49999
$ws.Range("L92").Value = 49999
$ws.Range("N92").Value = -53743

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = ""

$ws.Range("H97").Value = 504.64
$ws.Range("I97").Value = 541.3333
$ws.Range("J97").Value = 410.2857
$ws.Range("K97").Value = 541.3333
$ws.Range("L97").Value = 410.2857
$ws.Range("M97").Value = -45.33330000000001
$ws.Range("N97").Value = -1402.2857

$ws.Range("H98").Value = 38995
$ws.Range("J98").Value = 38995
$ws.Range("L98").Value = 38995
$ws.Range("N98").Value = -44985

$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""

$ws.Range("H105").Value = 52101
$ws.Range("J105").Value = 52101
$ws.Range("L105").Value = 52101
$ws.Range("N105").Value = -59089

$ws.Range("H107").Value = 167814.5
$ws.Range("I107").Value = 250346.75
$ws.Range("J107").Value = 2750
$ws.Range("K107").Value = 250346.75
$ws.Range("L107").Value = 2750
$ws.Range("M107").Value = -248426.75
$ws.Range("N107").Value = -6590

$ws.Range("H108").Value = 177777
$ws.Range("J108").Value = 177777
$ws.Range("L108").Value = 177777
$ws.Range("N108").Value = -185457

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = ""

$ws.Range("H113").Value = 4078.7917
$ws.Range("I113").Value = 3948.7368
$ws.Range("J113").Value = 4573
$ws.Range("K113").Value = 3948.7368
$ws.Range("L113").Value = 4573
$ws.Range("M113").Value = -1778.7368
$ws.Range("N113").Value = -8913

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1800.2667
$ws.Range("I107").Value = 2090.6667
$ws.Range("J107").Value = 1675.8096
$ws.Range("K107").Value = 6272.000100000001
$ws.Range("L107").Value = 5027.4288
$ws.Range("M107").Value = -4352.000100000001
$ws.Range("N107").Value = -8867.4288

$ws.Range("H126").Value = 2603.1333
$ws.Range("I126").Value = 2603.1333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7809.3999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5339.3999
$ws.Range("N126").Value = ""

$ws.Range("H136").Value = 3545.2083
$ws.Range("I136").Value = 2627.5833
$ws.Range("J136").Value = 4462.8335
$ws.Range("K136").Value = 7882.749899999999
$ws.Range("L136").Value = 13388.5005
$ws.Range("M136").Value = -5332.749899999999
$ws.Range("N136").Value = -18488.5005

$ws.Range("H139").Value = 173324.5
$ws.Range("J139").Value = 173324.5
$ws.Range("L139").Value = 173324.5
$ws.Range("N139").Value = -183604.5
